$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The teacher entry that used to live in row 25 (Reynaldo Abiog) was removed
# from the roster. Clear the row's values but keep the existing cell
# formatting for the remaining columns (A-D, F-I).
$ws.Range("A25:D25").ClearContents()
$ws.Range("F25:I25").ClearContents()

# Column E (the "Lumaniag, Lian, Batangas" address) had its formatting
# cleared entirely, not just its value, so the cell drops out of the sheet.
$ws.Range("E25").Clear()

# Drop the mailto hyperlink that pointed at the removed teacher's e-mail.
foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Range.Address($false, $false) -eq "F25") {
        $hl.Delete()
    }
}

# The row no longer wraps an address, so its custom height goes back to
# the sheet's default (same as AutoFit would give it).
$ws.Rows.Item(25).AutoFit()

# Leave the selection on the now-empty row.
$ws.Range("A25:I25").Select()
